$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing prices ---
$ws.Range("C6").Value = 2000
$ws.Range("C7").Value = 8000

# --- Add new product rows 50-57 (name / image URL / price) ---
$rows = @(
    @{ Row = 50; Name = "컵라면";         Url = "https://img.freepik.com/premium-psd/cup-instant-ramen-with-green-onions_787500-55485.jpg?ga=GA1.1.426326816.1750134496&semt=ais_hybrid&w=740"; Price = 1000 },
    @{ Row = 51; Name = "모기 기피제";     Url = "https://img.freepik.com/premium-photo/man-applies-insect-repellent-his-sons-hand-park-closeup-prevention-bites_472916-50499.jpg?ga=GA1.1.426326816.1750134496&semt=ais_hybrid&w=740"; Price = 3000 },
    @{ Row = 52; Name = "과자";           Url = "https://img.freepik.com/free-photo/crackers-biscuits_1339-1032.jpg?ga=GA1.1.426326816.1750134496&semt=ais_hybrid&w=740"; Price = 1000 },
    @{ Row = 53; Name = "탄산음료(콜라)"; Url = "https://img.freepik.com/premium-photo/red-aluminum-cans-with-water-droplets-white-background_167862-6023.jpg?ga=GA1.1.426326816.1750134496&semt=ais_hybrid&w=740"; Price = 1500 },
    @{ Row = 54; Name = "케이크";         Url = "https://encrypted-tbn0.gstatic.com/images?q=tbn:ANd9GcTF19OGL1O_Ky9P0hTPDxGxyVpvssPD1642jA&s"; Price = 20000 },
    @{ Row = 55; Name = "조각 케이크";    Url = "https://lucycato.co.kr/pds/space/125_1?1651652964"; Price = 5000 },
    @{ Row = 56; Name = "생일초";         Url = "https://encrypted-tbn0.gstatic.com/images?q=tbn:ANd9GcSG6PyRclVULM2jfcOCjWbC-PHTpgSwJOik0A&s"; Price = 1000 },
    @{ Row = 57; Name = "파티풍선 세트"; Url = "https://m.hiballoon.co.kr/web/product/big/202204/9cf4b19c088b4bf1c369133b24986082.jpg"; Price = 4000 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.Name
    $ws.Cells.Item($r, 3).Value = $item.Price

    $cell = $ws.Cells.Item($r, 2)
    $ws.Hyperlinks.Add($cell, $item.Url)
    # Hyperlinks.Add re-styles the cell; restore the sheet's standard hyperlink style
    $cell.Style = "하이퍼링크"
}

# --- Restore final selection state ---
$ws.Range("B63").Select()
